# Update benchmark: 2026-02-12 07:10:40 UTC
# Clear stale benchmark figures in column E/G for rows 2-14 (ZIRAAT / ISBANKASI columns)
# that are being refreshed / removed in this pass.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BENCHMARK")

$cellsToClear = @("G2", "G3", "G4", "G5", "E6", "G6", "G7", "G8", "G9", "G10", "G11", "G12", "E13", "E14", "G14")

foreach ($addr in $cellsToClear) {
    $ws.Range($addr).Value = ""
}
